$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.164.67'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.812.89'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.97%  '
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4609'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.70%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3746'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07385'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8638'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.58'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.821.18'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.652'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("E14").Value = '  +2.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07075'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.68'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008732'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.86'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.186.86'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.308'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.89'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.041.95'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.930'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.212'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.36%  '
$ws.Range("E28").Value = '  +1.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.258'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.83'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08891'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7705'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.170'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.510'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.912'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.27%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.125'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01958'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05232'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.268'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.382'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +21.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.913'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5279'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.81%  '
$ws.Range("E44").Value = '  +0.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.591'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5015'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '104.96'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.36%  '
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.668'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.89%  '
$ws.Range("E51").Value = '  +0.20%  '
